$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename existing site title: Skomakerfjellet -> Valsneset ---
$ws.Range("A19").Value = "Valsneset"

# --- Valsneset error breakdown (rows 20-37) ---
$ws.Range("A20").Value = "62 Feeding fault"
$ws.Range("B20").Value = 89
$ws.Range("A21").Value = "21 Cable twisted"
$ws.Range("B21").Value = 217
$ws.Range("A22").Value = "9 Generator Heating"
$ws.Range("B22").Value = 79
$ws.Range("A23").Value = "60 Mains failure"
$ws.Range("B23").Value = 85
$ws.Range("A24").Value = "25 Fault yaw inverter"
$ws.Range("B24").Value = 123
$ws.Range("A25").Value = "222:1 Turbine reset (power failure)"
$ws.Range("B25").Value = 27
$ws.Range("A26").Value = "222:1 Turbine reset (quit button)"
$ws.Range("B26").Value = 11773
$ws.Range("A27").Value = "220 Processor reset"
$ws.Range("B27").Value = 24
$ws.Range("A28").Value = "42 Pitch control error"
$ws.Range("B28").Value = 55
$ws.Range("A29").Value = "29 Anemonmeter Interface"
$ws.Range("B29").Value = 63
$ws.Range("A30").Value = "15 Turbine moist"
$ws.Range("B30").Value = 2
$ws.Range("A31").Value = "67 Overtemperature"
$ws.Range("B31").Value = 9
$ws.Range("A32").Value = "66 Fault rectifier"
$ws.Range("B32").Value = 1
$ws.Range("A33").Value = "72 Air gap monitoring"
$ws.Range("B33").Value = 3
$ws.Range("A34").Value = "48 Fault speed sensor"
$ws.Range("B34").Value = 8
$ws.Range("A35").Value = "31 Tower oscillation"
$ws.Range("B35").Value = 0
$ws.Range("A36").Value = "64 Overcurrent inverter"
$ws.Range("B36").Value = 0
$ws.Range("A37").Value = "70 Generator overtemperature"
$ws.Range("B37").Value = 43

# --- Ytre Vika title + error breakdown (rows 39-57) ---
$ws.Range("A39").Value = "Ytre Vika"
$ws.Range("A40").Value = "62 Feeding fault"
$ws.Range("B40").Value = 10
$ws.Range("A41").Value = "21 Cable twisted"
$ws.Range("B41").Value = 164
$ws.Range("A42").Value = "9 Generator Heating"
$ws.Range("B42").Value = 32
$ws.Range("A43").Value = "60 Mains failure"
$ws.Range("B43").Value = 0
$ws.Range("A44").Value = "25 Fault yaw inverter"
$ws.Range("B44").Value = 9
$ws.Range("A45").Value = "222:1 Turbine reset (power failure)"
$ws.Range("B45").Value = 17
$ws.Range("A46").Value = "222:1 Turbine reset (quit button)"
$ws.Range("B46").Value = 0
$ws.Range("A47").Value = "220 Processor reset"
$ws.Range("B47").Value = 16
$ws.Range("A48").Value = "42 Pitch control error"
$ws.Range("B48").Value = 9
$ws.Range("A49").Value = "29 Anemonmeter Interface"
$ws.Range("B49").Value = 5
$ws.Range("A50").Value = "15 Turbine moist"
$ws.Range("B50").Value = 0
$ws.Range("A51").Value = "67 Overtemperature"
$ws.Range("B51").Value = 15
$ws.Range("A52").Value = "66 Fault rectifier"
$ws.Range("B52").Value = 20
$ws.Range("A53").Value = "72 Air gap monitoring"
$ws.Range("B53").Value = 0
$ws.Range("A54").Value = "48 Fault speed sensor"
$ws.Range("B54").Value = 0
$ws.Range("A55").Value = "31 Tower oscillation"
$ws.Range("B55").Value = 0
$ws.Range("A56").Value = "64 Overcurrent inverter"
$ws.Range("B56").Value = 17
$ws.Range("A57").Value = "70 Generator overtemperature"
$ws.Range("B57").Value = 36

# --- Formatting: match the existing Bessaker blocks fonts ---
# Title rows (big bold) use the same look as A1/A19
$ws.Range("A39").Font.Bold = $true
$ws.Range("A39").Font.Size = 40
$ws.Range("A39").RowHeight = 50

# Name cells (column A): regular weight, size 22 (like A2:A17)
$ws.Range("A20:A37").Font.Bold = $false
$ws.Range("A20:A37").Font.Size = 22
$ws.Range("A40:A57").Font.Bold = $false
$ws.Range("A40:A57").Font.Size = 22

# Value cells (column B): bold, size 22 (like B2:B17)
$ws.Range("B20:B37").Font.Bold = $true
$ws.Range("B20:B37").Font.Size = 22
$ws.Range("B40:B57").Font.Bold = $true
$ws.Range("B40:B57").Font.Size = 22

# --- Leave the selection where the author left it ---
$ws.Range("I33").Select() | Out-Null
